$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 ("Emerald" / Lakeside Paddle Boats, Emerald Lake Park ...)
# This shifts all subsequent rows up by one, matching the diff which
# removes the Emerald entry and renumbers rows 29-54 to 28-53.
$ws.Rows.Item(28).Delete()
